$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows with newly received invoice/customer data ---

# Row 5: BOMBUS CONCEPTS LLC -> Prospect code, Last Invoice Date, Customer Number
$ws.Range("C5").Value = "023"
$ws.Range("D8").Copy($ws.Range("D5"))
$ws.Range("D5").Value = 45965
$ws.Range("E5").Value = "0008313"

# Row 18: BEHRMAV ENTERPRISES LLC -> Last Invoice Date, Customer Number
$ws.Range("D8").Copy($ws.Range("D18"))
$ws.Range("D18").Value = 45965
$ws.Range("E18").Value = "0008337"

# Row 24: ROSALIA LLC -> Last Invoice Date, Customer Number
$ws.Range("D8").Copy($ws.Range("D24"))
$ws.Range("D24").Value = 45965
$ws.Range("E24").Value = "0008350"

# --- Insert a brand new customer row (PETE AND PORKY BIG GAME) before row 26 ---
$ws.Rows.Item(26).Insert()

$ws.Range("A26").Value = "PETE AND PORKY BIG GAME"
$ws.Range("B26").Value = "Monroe, Michael D"
$ws.Range("C26").Value = "030"
$ws.Range("D8").Copy($ws.Range("D26"))
$ws.Range("D26").Value = 45966
$ws.Range("E26").Value = "0008359"
$ws.Rows.Item(26).RowHeight = 13.05
